$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from the last existing data row (171) down into the
# two new rows (172, 173) so fonts/number formats/borders match the rest of
# the table.
$ws.Range("A171:I171").Copy()
$ws.Range("A172:I173").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 172: LeetCode 2221 - Find Triangular Sum of an Array ---
$ws.Cells.Item(172, 1).Value2 = 2221
$ws.Cells.Item(172, 2).Value2 = "Find Triangular Sum of an Array"
$ws.Cells.Item(172, 3).Value2 = "#array #math #simulation #combination "
$ws.Cells.Item(172, 4).Value2 = "medium"
$ws.Cells.Item(172, 5).Value2 = 1
$ws.Cells.Item(172, 6).Value2 = 0
$ws.Cells.Item(172, 7).Value2 = 30
$ws.Cells.Item(172, 8).Value2 = 45930
$ws.Cells.Item(172, 9).Value2 = 45930
$ws.Rows.Item(172).RowHeight = 51

# --- Row 173: LeetCode 1121 - Divide Array Into Increasing Sequences ---
$ws.Cells.Item(173, 1).Value2 = 1121
$ws.Cells.Item(173, 2).Value2 = "Divide Array Into Increasing Sequences"
$ws.Cells.Item(173, 3).Value2 = "#array #counting"
$ws.Cells.Item(173, 4).Value2 = "hard"
$ws.Cells.Item(173, 5).Value2 = 0
$ws.Cells.Item(173, 6).Value2 = 1
$ws.Cells.Item(173, 7).Value2 = 20
$ws.Cells.Item(173, 8).Value2 = 45930
$ws.Cells.Item(173, 9).Value2 = 45930
$ws.Rows.Item(173).RowHeight = 34

# Restore view state: select C174 (the next empty row), matching where the
# user left off after entering the two new rows.
[void]$ws.Range("C174").Select()
